$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Absence")

# Replace the two absence columns' headers with the new "Nature absence" /
# "Nombre jour" headers and drop the old "Commentaires" column (E).
$ws1.Range("C1").Value = "Nature absence"
$ws1.Range("D1").Value = "Nombre jour"
$ws1.Columns.Item(5).Delete()

# The demo date value in B2 is no longer populated by this template.
$ws1.Range("B2").ClearContents()

# The "Nombre jour" column is right aligned.
$ws1.Range("D2").HorizontalAlignment = -4152

# Match the new column widths used by the updated template.
$ws1.Columns.Item(3).ColumnWidth = 15.166666666666666
$ws1.Columns.Item(4).ColumnWidth = 11.666666666666666

$ws1.Range("D2").Select() | Out-Null
